# Generate Report for Handoff
# Refreshes the GUID-based file names, content hash, and timestamps across
# the Overview / zh-cn / de-de sheets to reflect the latest handoff run.

$wb = $excel.ActiveWorkbook

$newGuid = "1442afcb-9d0e-4476-99c2-45f1560af881"
$newHash = "88c59249ca584fd795a79250abd1c6345bfec02d"

$hyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/09048dca770ec3c10889bb4606c8836e1fe9e905/e2e/1789a0f1-fed2-4424-ba61-47a9b82323ea.md"

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
$overview.Range("A2").Value = "$newGuid.md"
$overview.Range("B2").Value = "e2e\$newGuid.md"
$overview.Range("G2").Value = "2016-09-03 19:11:32"

$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add($overview.Range("B2"), $hyperlinkAddress, "", "", "e2e\$newGuid.md")
$overview.Range("B2").Font.Color = 15570276

# --- zh-cn sheet ---
$zhcn.Range("A2").Value = "$newGuid.md"
$zhcn.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$zhcn.Range("H2").Value = "2016-09-03 19:11:28"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$zhcn.Range("A2").Font.Color = 15570276

# --- de-de sheet ---
$dede.Range("A2").Value = "$newGuid.md"
$dede.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$dede.Range("H2").Value = "2016-09-03 19:11:32"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $hyperlinkAddress, "", "", "$newGuid.md")
$dede.Range("A2").Font.Color = 15570276
